$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values must stay text exactly as scraped (may look numeric,
# including trailing zeros like "1.000"), so force text format before writing,
# then clear the format back so no stray style is left behind.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.301.47'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.857.70'
$ws.Range('D3').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.31'
$ws.Range('D5').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4732'
$ws.Range('D7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2739'
$ws.Range('D8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06418'
$ws.Range('D9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.874.00'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07460'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.32'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.985'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '85.33'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6324'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.272.85'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('D17').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.25'
$ws.Range('D18').ClearFormats()
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.74'
$ws.Range('D19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007385'
$ws.Range('D20').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.101.27'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.987'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.972'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.280'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.85'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.97'
$ws.Range('D27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.886'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1046'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.399'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.150'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.924'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04928'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.162'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7236'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.701'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01869'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.647'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9141'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.969'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '105.87'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9998'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4103'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.550'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.107'
$ws.Range('D46').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1197'
$ws.Range('D48').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.699'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.41'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.405'
$ws.Range('D51').ClearFormats()

# Column E (Volume(1h)) values are padded percentage strings; plain text already.
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('E3').Value = '  -2.49%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('E5').Value = '  -2.00%  '
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('E8').Value = '  -3.41%  '
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('E10').Value = '  -6.54%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('E14').Value = '  -3.36%  '
$ws.Range('E15').Value = '  -5.47%  '
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('E19').Value = '  -4.49%  '
$ws.Range('E20').Value = '  -2.78%  '
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  -6.27%  '
$ws.Range('E24').Value = '  -4.16%  '
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('E26').Value = '  -1.27%  '
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('E28').Value = '  -3.81%  '
$ws.Range('E29').Value = '  +7.30%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  -4.97%  '
$ws.Range('E32').Value = '  -2.53%  '
$ws.Range('E33').Value = '  -2.79%  '
$ws.Range('E34').Value = '  -6.26%  '
$ws.Range('E35').Value = '  -4.21%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('E41').Value = '  -5.64%  '
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('E44').Value = '  -4.52%  '
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('E46').Value = '  -4.33%  '
$ws.Range('E47').Value = '  -5.86%  '
$ws.Range('E48').Value = '  -6.15%  '
$ws.Range('E49').Value = '  -3.77%  '
$ws.Range('E50').Value = '  -1.35%  '
$ws.Range('E51').Value = '  -5.10%  '
